# "Generate Report for Handback" — refresh the localization-status report:
#   - Status flips from "Ready for handoff" to "Handed back: in sync with en-US"
#     on the Overview summary sheet and on each language sheet.
#   - The stale "handback file is not latest" Error Detail is cleared now that
#     the handback is in sync.
#   - The Latest Handback DateTime timestamps are refreshed.
#   - The Status / Error Detail columns are widened/narrowed to fit the new text.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status column updates ---------------------------------------------
$wsOverview.Range("E2").Value = $newStatus   # zh-cn status
$wsOverview.Range("F2").Value = $newStatus   # de-de status
$wsZhCn.Range("C2").Value     = $newStatus
$wsDeDe.Range("C2").Value     = $newStatus

# --- Error Detail cleared (handback is now in sync, no more errors) -----
$wsZhCn.Range("P2").Value = ""
$wsDeDe.Range("P2").Value = ""

# --- Latest Handback DateTime refreshed ---------------------------------
$wsZhCn.Range("K2").Value = "2016-09-05 15:02:34"
$wsDeDe.Range("K2").Value = "2016-09-05 15:02:42"

# --- Column widths: autofit-style resize for the edited columns --------
# Status columns (now longer text) widen; Error Detail (now empty) narrows.
$wsOverview.Columns.Item(5).ColumnWidth = 29.15   # E: zh-cn status
$wsOverview.Columns.Item(6).ColumnWidth = 29.15   # F: de-de status

$wsZhCn.Columns.Item(3).ColumnWidth  = 29.15      # C: Status
$wsZhCn.Columns.Item(16).ColumnWidth = 12.85      # P: Error Detail

$wsDeDe.Columns.Item(3).ColumnWidth  = 29.15      # C: Status
$wsDeDe.Columns.Item(16).ColumnWidth = 12.85      # P: Error Detail
